$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert 18 new rows (206-223), each copied from the last existing
# data row (205) so they inherit its cell style (s="1") on all 5 columns.
for ($i = 0; $i -lt 18; $i++) {
    $destRow = 206 + $i
    $ws.Rows.Item(205).Copy()
    $ws.Rows.Item($destRow).Insert(-4121)
}

# Step 2: fill in the values for each of the new rows.
# Row 206
$ws.Cells.Item(206,1).Value = "AppointmentRead"
$ws.Cells.Item(206,2).Value = "end"
$ws.Cells.Item(206,4).Value = "outlook-other-item-apis-get-end-read"
$ws.Cells.Item(206,5).Value = "get"

# Row 207
$ws.Cells.Item(207,1).Value = "MessageRead"
$ws.Cells.Item(207,2).Value = "end"
$ws.Cells.Item(207,4).Value = "outlook-other-item-apis-get-end-read"
$ws.Cells.Item(207,5).Value = "get"

# Row 208
$ws.Cells.Item(208,1).Value = "AppointmentCompose"
$ws.Cells.Item(208,2).Value = "end"
$ws.Cells.Item(208,4).Value = "outlook-other-item-apis-get-set-end-appointment-organizer"
$ws.Cells.Item(208,5).Value = "get"

# Row 209
$ws.Cells.Item(209,1).Value = "AppointmentCompose"
$ws.Cells.Item(209,2).Value = "end"
$ws.Cells.Item(209,4).Value = "outlook-other-item-apis-get-set-end-appointment-organizer"
$ws.Cells.Item(209,5).Value = "set"

# Row 210
$ws.Cells.Item(210,1).Value = "Time"
$ws.Cells.Item(210,2).Value = "setAsync"
$ws.Cells.Item(210,3).Value = 2
$ws.Cells.Item(210,4).Value = "outlook-other-item-apis-get-set-end-appointment-organizer"
$ws.Cells.Item(210,5).Value = "set"

# Row 211
$ws.Cells.Item(211,1).Value = "AppointmentRead"
$ws.Cells.Item(211,2).Value = "location"
$ws.Cells.Item(211,4).Value = "outlook-other-item-apis-get-location-read"
$ws.Cells.Item(211,5).Value = "get"

# Row 212
$ws.Cells.Item(212,1).Value = "MessageRead"
$ws.Cells.Item(212,2).Value = "location"
$ws.Cells.Item(212,4).Value = "outlook-other-item-apis-get-location-read"
$ws.Cells.Item(212,5).Value = "get"

# Row 213
$ws.Cells.Item(213,1).Value = "AppointmentCompose"
$ws.Cells.Item(213,2).Value = "location"
$ws.Cells.Item(213,4).Value = "outlook-other-item-apis-get-set-location-appointment-organizer"
$ws.Cells.Item(213,5).Value = "get"

# Row 214
$ws.Cells.Item(214,1).Value = "Location"
$ws.Cells.Item(214,2).Value = "getAsync"
$ws.Cells.Item(214,3).Value = 2
$ws.Cells.Item(214,4).Value = "outlook-other-item-apis-get-set-location-appointment-organizer"
$ws.Cells.Item(214,5).Value = "get"

# Row 215
$ws.Cells.Item(215,1).Value = "AppointmentCompose"
$ws.Cells.Item(215,2).Value = "location"
$ws.Cells.Item(215,4).Value = "outlook-other-item-apis-get-set-location-appointment-organizer"
$ws.Cells.Item(215,5).Value = "set"

# Row 216
$ws.Cells.Item(216,1).Value = "Location"
$ws.Cells.Item(216,2).Value = "setAsync"
$ws.Cells.Item(216,3).Value = 2
$ws.Cells.Item(216,4).Value = "outlook-other-item-apis-get-set-location-appointment-organizer"
$ws.Cells.Item(216,5).Value = "set"

# Row 217
$ws.Cells.Item(217,1).Value = "AppointmentCompose"
$ws.Cells.Item(217,2).Value = "enhancedLocation"
$ws.Cells.Item(217,4).Value = "outlook-other-item-apis-get-add-remove-enhancedlocation-appointment"
$ws.Cells.Item(217,5).Value = "get"

# Row 218
$ws.Cells.Item(218,1).Value = "AppointmentRead"
$ws.Cells.Item(218,2).Value = "enhancedLocation"
$ws.Cells.Item(218,4).Value = "outlook-other-item-apis-get-add-remove-enhancedlocation-appointment"
$ws.Cells.Item(218,5).Value = "get"

# Row 219
$ws.Cells.Item(219,1).Value = "EnhancedLocation"
$ws.Cells.Item(219,2).Value = "getAsync"
$ws.Cells.Item(219,3).Value = 2
$ws.Cells.Item(219,4).Value = "outlook-other-item-apis-get-add-remove-enhancedlocation-appointment"
$ws.Cells.Item(219,5).Value = "get"

# Row 220
$ws.Cells.Item(220,1).Value = "AppointmentCompose"
$ws.Cells.Item(220,2).Value = "enhancedLocation"
$ws.Cells.Item(220,4).Value = "outlook-other-item-apis-get-add-remove-enhancedlocation-appointment"
$ws.Cells.Item(220,5).Value = "add"

# Row 221
$ws.Cells.Item(221,1).Value = "EnhancedLocation"
$ws.Cells.Item(221,2).Value = "addAsync"
$ws.Cells.Item(221,3).Value = 2
$ws.Cells.Item(221,4).Value = "outlook-other-item-apis-get-add-remove-enhancedlocation-appointment"
$ws.Cells.Item(221,5).Value = "add"

# Row 222
$ws.Cells.Item(222,1).Value = "AppointmentCompose"
$ws.Cells.Item(222,2).Value = "enhancedLocation"
$ws.Cells.Item(222,4).Value = "outlook-other-item-apis-get-add-remove-enhancedlocation-appointment"
$ws.Cells.Item(222,5).Value = "remove"

# Row 223
$ws.Cells.Item(223,1).Value = "EnhancedLocation"
$ws.Cells.Item(223,2).Value = "removeAsync"
$ws.Cells.Item(223,3).Value = 2
$ws.Cells.Item(223,4).Value = "outlook-other-item-apis-get-add-remove-enhancedlocation-appointment"
$ws.Cells.Item(223,5).Value = "remove"

# Step 3: row 220 is special-cased in the source data - columns A, B and E
# carry no explicit style (unlike every other new row) and column C has no
# cell at all, while column D keeps the inherited style.
$ws.Range("A220:B220").ClearFormats()
$ws.Range("E220").ClearFormats()
$ws.Cells.Item(220,3).ClearFormats()
$ws.Cells.Item(220,3).ClearContents()

# Re-apply the values on row 220 that ClearFormats() may have left untouched
# (ClearFormats does not remove values, only formatting, so values already set
# above remain intact; this is just defensive re-assertion).
$ws.Cells.Item(220,1).Value = "AppointmentCompose"
$ws.Cells.Item(220,2).Value = "enhancedLocation"
$ws.Cells.Item(220,4).Value = "outlook-other-item-apis-get-add-remove-enhancedlocation-appointment"
$ws.Cells.Item(220,5).Value = "add"

# Step 4: resize the "Snippets" table and its autofilter to cover the new rows.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E223"))

# Step 5: update the selection / active cell to the first cell of the next
# (still empty) row, matching the end-user state captured in the workbook.
$ws.Range("A224").Select()
